# Update cryptocurrency price/volume figures (and one coin swap in row 51)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'30.213.13"
$ws.Cells.Item(2, 5).Value = "'  -0.58%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'2.064.55"
$ws.Cells.Item(3, 5).Value = "'  +2.97%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "'  +0.04%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'325.97"
$ws.Cells.Item(5, 5).Value = "'  +0.46%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "'  +0.18%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.5157"
$ws.Cells.Item(7, 5).Value = "'  +1.25%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "'  +3.88%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.08702"
$ws.Cells.Item(9, 5).Value = "'  -0.56%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'45.43"
$ws.Cells.Item(10, 5).Value = "'  +5.53%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "'  +1.22%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "'  -2.35%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'2.061.16"
$ws.Cells.Item(13, 5).Value = "'  +2.82%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'6.604"
$ws.Cells.Item(14, 5).Value = "'  +0.60%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'7.599"
$ws.Cells.Item(15, 5).Value = "'  +1.93%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'1.005"
$ws.Cells.Item(16, 5).Value = "'  +0.33%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'94.47"
$ws.Cells.Item(17, 5).Value = "'  +0.24%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.00001113"
$ws.Cells.Item(18, 5).Value = "'  -0.14%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'0.06604"
$ws.Cells.Item(19, 5).Value = "'  +1.37%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'18.63"
$ws.Cells.Item(20, 5).Value = "'  -1.40%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'1.002"
$ws.Cells.Item(21, 5).Value = "'  +0.08%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.188"
$ws.Cells.Item(22, 5).Value = "'  -0.37%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'30.222.15"
$ws.Cells.Item(23, 5).Value = "'  -0.72%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "'  +1.74%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.273"
$ws.Cells.Item(25, 5).Value = "'  +2.18%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'2.305.39"
$ws.Cells.Item(26, 5).Value = "'  +3.13%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'22.02"
$ws.Cells.Item(27, 5).Value = "'  -1.19%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'162.47"
$ws.Cells.Item(28, 5).Value = "'  -0.20%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'2.481"
$ws.Cells.Item(29, 5).Value = "'  +3.05%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'129.98"

# Row 31
$ws.Cells.Item(31, 4).Value = "'1.163"
$ws.Cells.Item(31, 5).Value = "'  +2.67%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.1058"
$ws.Cells.Item(32, 5).Value = "'  +0.66%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'6.037"
$ws.Cells.Item(33, 5).Value = "'  -0.52%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'3.835"
$ws.Cells.Item(34, 5).Value = "'  +0.11%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.489"
$ws.Cells.Item(35, 5).Value = "'  +10.56%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.02541"
$ws.Cells.Item(36, 5).Value = "'  +0.87%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'9.526"
$ws.Cells.Item(37, 5).Value = "'  +5.57%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'5.385"
$ws.Cells.Item(38, 5).Value = "'  -0.72%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.06550"
$ws.Cells.Item(39, 5).Value = "'  -0.58%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'12.43"
$ws.Cells.Item(40, 5).Value = "'  +0.09%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.2216"
$ws.Cells.Item(41, 5).Value = "'  +1.03%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'0.6592"
$ws.Cells.Item(42, 5).Value = "'  -0.63%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'1.229"
$ws.Cells.Item(43, 5).Value = "'  +0.31%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "'  +0.24%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'13.89"
$ws.Cells.Item(45, 5).Value = "'  +2.92%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.6246"
$ws.Cells.Item(46, 5).Value = "'  +1.45%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'2.172"
$ws.Cells.Item(47, 5).Value = "'  -0.56%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'3.590"
$ws.Cells.Item(48, 5).Value = "'  -1.99%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'1.229"
$ws.Cells.Item(49, 5).Value = "'  -2.58%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'80.77"
$ws.Cells.Item(50, 5).Value = "'  +0.69%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "'Quant"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(51, 4).Value = "'119.63"
$ws.Cells.Item(51, 5).Value = "'  -3.70%  "
